$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Anwesend" attendance cells, copying the formatting from
# the neighboring cell in the same row that already carries the fill.
$ws.Range("L5").Copy($ws.Range("M5"))
$ws.Range("M5").Value = "Anwesend"

$ws.Range("L6").Copy($ws.Range("M6"))
$ws.Range("M6").Value = "Anwesend"

$ws.Range("L7").Copy($ws.Range("M7"))
$ws.Range("M7").Value = "Anwesend"

$ws.Range("K8").Copy($ws.Range("L8"))
$ws.Range("L8").Value = "Anwesend"
$ws.Range("K8").Copy($ws.Range("M8"))
$ws.Range("M8").Value = "Anwesend"

$ws.Range("L9").Copy($ws.Range("M9"))
$ws.Range("M9").Value = "Anwesend"

# Update the visible sheet view: scroll and selection moved.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("N5").Select()
